# Onward Technologies Ltd_Cleaned_Data.xlsx
# "new added data and sankey"
#
# The "Quarterly" sheet gains a new "Exceptional items" column between the
# existing "P/l before exceptional items & tax" column (K/K) and
# "P/l before tax" column (old L). We insert a blank column at L, which
# shifts everything from L..T to M..U, then populate the two header rows
# for the freshly inserted column:
#   - Row 1 (canonical/lowercase labels): "Exceptional items"
#   - Row 2 (display labels)            : "Exceptional Items"
# The data rows (3-47) are left blank in the new column, matching the
# source data (no historical "exceptional items" figures were available).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new blank column before column L; shifts L:T -> M:U.
$ws.Columns("L").Insert()

# Header row 1 (bold/bordered style already carried over from the Insert).
$ws.Range("L1").Value = "Exceptional items"

# Header row 2 (plain style, matches its row-2 siblings).
$ws.Range("L2").Value = "Exceptional Items"
